$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -13.16
$ws.Range("B7").Value = 5.910000000000001
$ws.Range("A8").Value = -22.084
$ws.Range("A10").Value = -21.657
$ws.Range("E10").Value = 16.331
$ws.Range("A12").Value = -21.475
$ws.Range("E12").Value = 17.587
$ws.Range("E13").Value = 16.542
$ws.Range("E14").Value = 16.816
$ws.Range("B15").Value = 5.072000000000001
$ws.Range("A18").Value = -21.324
$ws.Range("B18").Value = 7.621
$ws.Range("C18").Value = -11.793
$ws.Range("C19").Value = -11.706
$ws.Range("B20").Value = 6.956999999999999
$ws.Range("C27").Value = -13.447
$ws.Range("B29").Value = 5.42
$ws.Range("E29").Value = 17.233
$ws.Range("B30").Value = 6.208
$ws.Range("B31").Value = 5.522
$ws.Range("C31").Value = -13.323
$ws.Range("E32").Value = 16.674
$ws.Range("E35").Value = 16.384
$ws.Range("A37").Value = -20.287
$ws.Range("C38").Value = -13.187
$ws.Range("B40").Value = 8.962
$ws.Range("C42").Value = -12.533
$ws.Range("E43").Value = 16.999
$ws.Range("C44").Value = -12.777
$ws.Range("C47").Value = -12.661
$ws.Range("E48").Value = 17.221
$ws.Range("E49").Value = 16.359
$ws.Range("B50").Value = 5.446
$ws.Range("E50").Value = 16.448
$ws.Range("A55").Value = -21.793
$ws.Range("E56").Value = 16.247
$ws.Range("C58").Value = -13.09
$ws.Range("C65").Value = -12.282
$ws.Range("A68").Value = -21.681
$ws.Range("B68").Value = 5.348000000000001
$ws.Range("E69").Value = 17.4
$ws.Range("C73").Value = -12.57
$ws.Range("B76").Value = 6.689
$ws.Range("A77").Value = -20.252
$ws.Range("A78").Value = -19.842
$ws.Range("A81").Value = -21.73
$ws.Range("E81").Value = 16.761
$ws.Range("A82").Value = -21.795
$ws.Range("B87").Value = 5.508999999999999
$ws.Range("B88").Value = 5.816000000000001
$ws.Range("C90").Value = -13.321
$ws.Range("E92").Value = 17.843
$ws.Range("C94").Value = -10.317
$ws.Range("C95").Value = -11.283
$ws.Range("B96").Value = 6.88
$ws.Range("B98").Value = 5.49
$ws.Range("B101").Value = 7.553999999999999
$ws.Range("C101").Value = -12.709
$ws.Range("B102").Value = 7.375999999999999
